$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the value from C7 (the Day 4 row had an extra, stray PR value of 200)
$ws.Range("C7").ClearContents()

# Move active selection to C7 to match the saved view state
$ws.Range("C7").Select()
